$d = $word.ActiveDocument

# Remove the "Scenario" through "Guidelines" through "Project Structure
# Overview" through "Requirements" heading sections (paragraphs 2-23),
# collapsing the document down to just the title and the requirements
# list that used to follow the "Requirements" heading.
$start = $d.Paragraphs.Item(2).Range.Start
$end = $d.Paragraphs.Item(23).Range.End
$r = $d.Range($start, $end)
$r.Delete()

# Retitle the document from "Fullstack Code Challenge" to "Requirments".
$d.Paragraphs.Item(1).Range.Text = "Requirments"
